$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("out_vars")

# Add a new row of data for 2020-08-02 (row 64), following the existing
# Raw/Clean SSA data layout: Fecha, Confirmados, Negativos, Sospechosos, Defunciones, Porcentaje hospitalizados
$row = 64
$ws.Cells.Item($row, 1).Value = "'2020-08-02"
$ws.Cells.Item($row, 2).Value = 439046
$ws.Cells.Item($row, 3).Value = 483333
$ws.Cells.Item($row, 4).Value = 83119
$ws.Cells.Item($row, 5).Value = 47746
$ws.Cells.Item($row, 6).Value = 26.92
